$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.098496255502688257
$ws.Range("A2").Value = -0.0059999999513635771
$ws.Range("A3").Value = -0.0039999999608832937
$ws.Range("A4").Value = -0.0079999999273141498
$ws.Range("A5").Value = 0.070455891250155389
$ws.Range("A6").Value = -0.0019999999611215458
$ws.Range("A7").Value = -0.0099999998975563997
$ws.Range("A8").Value = -0.0099999998957551739
$ws.Range("A9").Value = -0.0019999999587079209
$ws.Range("A10").Value = -0.0019999999583433237
$ws.Range("A11").Value = -0.0029999999504202179
$ws.Range("A12").Value = -0.003499999946503074
$ws.Range("A13").Value = -0.0034999999474845112
$ws.Range("A14").Value = -0.0079999999124380494
$ws.Range("A15").Value = 0.020196047175708998
$ws.Range("A16").Value = -0.0019999999607267505
$ws.Range("A17").Value = -0.0019999999599980001
$ws.Range("A18").Value = -0.003999999943889776
$ws.Range("A19").Value = -0.0039999999672732933
$ws.Range("A20").Value = 0.014873820281875894
$ws.Range("A21").Value = -0.0039999999677879927
$ws.Range("A22").Value = -0.0039999999675472964
$ws.Range("A23").Value = -0.0049999999500904835
$ws.Range("A24").Value = -0.019999999825960124
$ws.Range("A25").Value = -0.019999999823443027
$ws.Range("A26").Value = -0.0024999999554751184
$ws.Range("A27").Value = -0.0024999999531871708
$ws.Range("A28").Value = -0.0019999999461024487
$ws.Range("A29").Value = -0.0069999998986833845
$ws.Range("A30").Value = -0.033796738400227344
$ws.Range("A31").Value = -0.0069999998912138039
$ws.Range("A32").Value = -0.0099999998668600654
$ws.Range("A33").Value = -0.0039999999145337028
